$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.749.03"
$ws.Range("E2").Value = "  -3.26%  "
$ws.Range("D3").Value = "2.912.32"
$ws.Range("E3").Value = "  -3.90%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.57"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.40"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.29%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.503"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.93%  "
$ws.Range("D9").Value = "2.910.67"
$ws.Range("E9").Value = "  -3.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.77"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.94%  "
$ws.Range("E11").Value = "  -4.61%  "
$ws.Range("E12").Value = "  -3.03%  "
$ws.Range("E13").Value = "  -3.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.65"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.128"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("D16").Value = "3.395.93"
$ws.Range("E16").Value = "  -3.86%  "
$ws.Range("D17").Value = "60.749.98"
$ws.Range("E17").Value = "  -3.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.74"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.04%  "
$ws.Range("D19").Value = "2.912.47"
$ws.Range("E19").Value = "  -4.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "427.29"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -5.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.60"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.674"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.10"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.49"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.84"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.19"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.63%  "
$ws.Range("E27").Value = "  -2.39%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.17"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.92%  "
$ws.Range("E31").Value = "  -3.57%  "
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.46"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.85%  "
$ws.Range("E34").Value = "  -4.36%  "
$ws.Range("D35").Value = "0.0₃0856"
$ws.Range("E35").Value = "  +1.10%  "
$ws.Range("E36").Value = "  -2.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.63"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.99"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.56"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.125"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.00"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.65"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.288"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.02"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "375.87"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.46%  "
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("D47").Value = "2.664.88"
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.85"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.35"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("E51").Value = "  -1.77%  "
